$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.26539866666667
$ws.Range("H2").Value = 39.796196
$ws.Range("I2").Value = 0.164744212542501
$ws.Range("J2").Value = 0.164744212542501
$ws.Range("M2").Value = 5.256051666666667
$ws.Range("N2").Value = 15.768155
$ws.Range("O2").Value = 0.2751978571638378
$ws.Range("P2").Value = 0.2751978571638379
$ws.Range("Q2").Value = 69.72362077093112
$ws.Range("R2").Value = 627.5125869383801
$ws.Range("S2").Value = 0.04533725427184013
$ws.Range("T2").Value = 0.04533725427184013

$ws.Range("G3").Value = 13.26539866666667
$ws.Range("H3").Value = 39.796196
$ws.Range("I3").Value = 0.164744212542501
$ws.Range("J3").Value = 0.164744212542501
$ws.Range("O3").Value = 0.1142283101567343
$ws.Range("P3").Value = 0.1142283101567343
$ws.Range("Q3").Value = 28.94067366931156
$ws.Range("R3").Value = 260.466063023804
$ws.Range("S3").Value = 0.01881845300683177
$ws.Range("T3").Value = 0.01881845300683177

$ws.Range("G4").Value = 13.26539866666667
$ws.Range("H4").Value = 39.796196
$ws.Range("I4").Value = 0.164744212542501
$ws.Range("J4").Value = 0.164744212542501
$ws.Range("M4").Value = 11.66145566666667
$ws.Range("N4").Value = 34.98436700000001
$ws.Range("O4").Value = 0.6105738326794278
$ws.Range("P4").Value = 0.6105738326794279
$ws.Range("Q4").Value = 154.6938584519925
$ws.Range("R4").Value = 1392.244726067932
$ws.Range("S4").Value = 0.1005885052638291
$ws.Range("T4").Value = 0.1005885052638291

$ws.Range("I5").Value = 0.6678031736949381
$ws.Range("J5").Value = 0.6678031736949381
$ws.Range("M5").Value = 5.256051666666667
$ws.Range("N5").Value = 15.768155
$ws.Range("O5").Value = 0.2751978571638378
$ws.Range("P5").Value = 0.2751978571638379
$ws.Range("Q5").Value = 282.6299905395344
$ws.Range("R5").Value = 2543.66991485581
$ws.Range("S5").Value = 0.1837780024080572
$ws.Range("T5").Value = 0.1837780024080572

$ws.Range("I6").Value = 0.6678031736949381
$ws.Range("J6").Value = 0.6678031736949381
$ws.Range("O6").Value = 0.1142283101567343
$ws.Range("P6").Value = 0.1142283101567343
$ws.Range("S6").Value = 0.07628202804847692
$ws.Range("T6").Value = 0.07628202804847693

$ws.Range("I7").Value = 0.6678031736949381
$ws.Range("J7").Value = 0.6678031736949381
$ws.Range("M7").Value = 11.66145566666667
$ws.Range("N7").Value = 34.98436700000001
$ws.Range("O7").Value = 0.6105738326794278
$ws.Range("P7").Value = 0.6105738326794279
$ws.Range("Q7").Value = 627.0633003190039
$ws.Range("R7").Value = 5643.569702871035
$ws.Range("S7").Value = 0.407743143238404
$ws.Range("T7").Value = 0.407743143238404

$ws.Range("G8").Value = 13.48348233333333
$ws.Range("H8").Value = 40.450447
$ws.Range("I8").Value = 0.1674526137625609
$ws.Range("J8").Value = 0.1674526137625609
$ws.Range("M8").Value = 5.256051666666667
$ws.Range("N8").Value = 15.768155
$ws.Range("O8").Value = 0.2751978571638378
$ws.Range("P8").Value = 0.2751978571638379
$ws.Range("Q8").Value = 70.86987979058722
$ws.Range("R8").Value = 637.828918115285
$ws.Range("S8").Value = 0.04608260048394054
$ws.Range("T8").Value = 0.04608260048394054

$ws.Range("G9").Value = 13.48348233333333
$ws.Range("H9").Value = 40.450447
$ws.Range("I9").Value = 0.1674526137625609
$ws.Range("J9").Value = 0.1674526137625609
$ws.Range("O9").Value = 0.1142283101567343
$ws.Range("P9").Value = 0.1142283101567343
$ws.Range("Q9").Value = 29.41645946272811
$ws.Range("R9").Value = 264.748135164553
$ws.Range("S9").Value = 0.01912782910142565
$ws.Range("T9").Value = 0.01912782910142564

$ws.Range("G10").Value = 13.48348233333333
$ws.Range("H10").Value = 40.450447
$ws.Range("I10").Value = 0.1674526137625609
$ws.Range("J10").Value = 0.1674526137625609
$ws.Range("M10").Value = 11.66145566666667
$ws.Range("N10").Value = 34.98436700000001
$ws.Range("O10").Value = 0.6105738326794278
$ws.Range("P10").Value = 0.6105738326794279
$ws.Range("Q10").Value = 157.2370314624499
$ws.Range("R10").Value = 1415.133283162049
$ws.Range("S10").Value = 0.1022421841771947
$ws.Range("T10").Value = 0.1022421841771947
